$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$vals = New-Object 'object[,]' 18,16
$vals[0,0] = 2
$vals[0,1] = 0.6666666666666666
$vals[0,2] = 0.178715
$vals[0,3] = 0.536145
$vals[0,4] = 0.09904930989061336
$vals[0,5] = 0.09904930989061336
$vals[0,6] = 3
$vals[0,7] = 1
$vals[0,8] = 6.101885666666667
$vals[0,9] = 18.305657
$vals[0,10] = 0.1093737608697887
$vals[0,11] = 0.1093737608697887
$vals[0,12] = 1.090498496918333
$vals[0,13] = 9.814486472265
$vals[0,14] = 0.01083339553429355
$vals[0,15] = 0.01083339553429355
$vals[1,0] = 2
$vals[1,1] = 0.6666666666666666
$vals[1,2] = 0.178715
$vals[1,3] = 0.536145
$vals[1,4] = 0.09904930989061336
$vals[1,5] = 0.09904930989061336
$vals[1,6] = 3
$vals[1,7] = 1
$vals[1,8] = 29.178763
$vals[1,9] = 87.53628900000001
$vals[1,10] = 0.5230171820937495
$vals[1,11] = 0.5230171820937495
$vals[1,12] = 5.214682629545
$vals[1,13] = 46.932143665905
$vals[1,14] = 0.05180449094731915
$vals[1,15] = 0.05180449094731915
$vals[2,0] = 2
$vals[2,1] = 0.6666666666666666
$vals[2,2] = 0.178715
$vals[2,3] = 0.536145
$vals[2,4] = 0.09904930989061336
$vals[2,5] = 0.09904930989061336
$vals[2,6] = 2
$vals[2,7] = 0.6666666666666666
$vals[2,8] = 0.146644
$vals[2,9] = 0.439932
$vals[2,10] = 0.002628532664354407
$vals[2,11] = 0.002628532664354407
$vals[2,12] = 0.02620748246
$vals[2,13] = 0.23586734214
$vals[2,14] = 0.0002603543464292393
$vals[2,15] = 0.0002603543464292393
$vals[3,0] = 2
$vals[3,1] = 0.6666666666666666
$vals[3,2] = 0.178715
$vals[3,3] = 0.536145
$vals[3,4] = 0.09904930989061336
$vals[3,5] = 0.09904930989061336
$vals[3,6] = 3
$vals[3,7] = 1
$vals[3,8] = 15.02284966666667
$vals[3,9] = 45.068549
$vals[3,10] = 0.2692783275177917
$vals[3,11] = 0.2692783275177917
$vals[3,12] = 2.684808578178333
$vals[3,13] = 24.163277203605
$vals[3,14] = 0.02667183250913583
$vals[3,15] = 0.02667183250913583
$vals[4,0] = 2
$vals[4,1] = 0.6666666666666666
$vals[4,2] = 0.178715
$vals[4,3] = 0.536145
$vals[4,4] = 0.09904930989061336
$vals[4,5] = 0.09904930989061336
$vals[4,6] = 3
$vals[4,7] = 1
$vals[4,8] = 5.288900666666667
$vals[4,9] = 15.866702
$vals[4,10] = 0.09480134312252211
$vals[4,11] = 0.09480134312252211
$vals[4,12] = 0.9452058826433333
$vals[4,13] = 8.50685294379
$vals[4,14] = 0.009390007612989061
$vals[4,15] = 0.009390007612989061
$vals[5,0] = 2
$vals[5,1] = 0.6666666666666666
$vals[5,2] = 0.178715
$vals[5,3] = 0.536145
$vals[5,4] = 0.09904930989061336
$vals[5,5] = 0.09904930989061336
$vals[5,6] = 1
$vals[5,7] = 0.3333333333333333
$vals[5,8] = 0.050258
$vals[5,9] = 0.150774
$vals[5,10] = 0.0009008537317934847
$vals[5,11] = 0.0009008537317934848
$vals[5,12] = 0.00898185847
$vals[5,13] = 0.08083672622999999
$vals[5,14] = 0.00008922894044652836
$vals[5,15] = 0.00008922894044652837
$vals[6,0] = 3
$vals[6,1] = 1
$vals[6,2] = 1.450498333333333
$vals[6,3] = 4.351495
$vals[6,4] = 0.8039104659046613
$vals[6,5] = 0.8039104659046612
$vals[6,6] = 3
$vals[6,7] = 1
$vals[6,8] = 6.101885666666667
$vals[6,9] = 18.305657
$vals[6,10] = 0.1093737608697887
$vals[6,11] = 0.1093737608697887
$vals[6,12] = 8.850774989690555
$vals[6,13] = 79.656974907215
$vals[6,14] = 0.08792671105857688
$vals[6,15] = 0.08792671105857687
$vals[7,0] = 3
$vals[7,1] = 1
$vals[7,2] = 1.450498333333333
$vals[7,3] = 4.351495
$vals[7,4] = 0.8039104659046613
$vals[7,5] = 0.8039104659046612
$vals[7,6] = 3
$vals[7,7] = 1
$vals[7,8] = 29.178763
$vals[7,9] = 87.53628900000001
$vals[7,10] = 0.5230171820937495
$vals[7,11] = 0.5230171820937495
$vals[7,12] = 42.32374710022834
$vals[7,13] = 380.913723902055
$vals[7,14] = 0.4204589865331292
$vals[7,15] = 0.4204589865331292
$vals[8,0] = 3
$vals[8,1] = 1
$vals[8,2] = 1.450498333333333
$vals[8,3] = 4.351495
$vals[8,4] = 0.8039104659046613
$vals[8,5] = 0.8039104659046612
$vals[8,6] = 2
$vals[8,7] = 0.6666666666666666
$vals[8,8] = 0.146644
$vals[8,9] = 0.439932
$vals[8,10] = 0.002628532664354407
$vals[8,11] = 0.002628532664354407
$vals[8,12] = 0.2127068775933333
$vals[8,13] = 1.91436189834
$vals[8,14] = 0.002113104918846772
$vals[8,15] = 0.002113104918846772
$vals[9,0] = 3
$vals[9,1] = 1
$vals[9,2] = 1.450498333333333
$vals[9,3] = 4.351495
$vals[9,4] = 0.8039104659046613
$vals[9,5] = 0.8039104659046612
$vals[9,6] = 3
$vals[9,7] = 1
$vals[9,8] = 15.02284966666667
$vals[9,9] = 45.068549
$vals[9,10] = 0.2692783275177917
$vals[9,11] = 0.2692783275177917
$vals[9,12] = 21.79061840341722
$vals[9,13] = 196.115565630755
$vals[9,14] = 0.2164756657328559
$vals[9,15] = 0.2164756657328559
$vals[10,0] = 3
$vals[10,1] = 1
$vals[10,2] = 1.450498333333333
$vals[10,3] = 4.351495
$vals[10,4] = 0.8039104659046613
$vals[10,5] = 0.8039104659046612
$vals[10,6] = 3
$vals[10,7] = 1
$vals[10,8] = 5.288900666666667
$vals[10,9] = 15.866702
$vals[10,10] = 0.09480134312252211
$vals[10,11] = 0.09480134312252211
$vals[10,12] = 7.671541602165555
$vals[10,13] = 69.04387441949
$vals[10,14] = 0.07621179191801442
$vals[10,15] = 0.0762117919180144
$vals[11,0] = 3
$vals[11,1] = 1
$vals[11,2] = 1.450498333333333
$vals[11,3] = 4.351495
$vals[11,4] = 0.8039104659046613
$vals[11,5] = 0.8039104659046612
$vals[11,6] = 1
$vals[11,7] = 0.3333333333333333
$vals[11,8] = 0.050258
$vals[11,9] = 0.150774
$vals[11,10] = 0.0009008537317934847
$vals[11,11] = 0.0009008537317934848
$vals[11,12] = 0.07289914523666666
$vals[11,13] = 0.65609230713
$vals[11,14] = 0.0007242057432380531
$vals[11,15] = 0.000724205743238053
$vals[12,0] = 1
$vals[12,1] = 0.3333333333333333
$vals[12,2] = 0.17509
$vals[12,3] = 0.52527
$vals[12,4] = 0.09704022420472538
$vals[12,5] = 0.09704022420472537
$vals[12,6] = 3
$vals[12,7] = 1
$vals[12,8] = 6.101885666666667
$vals[12,9] = 18.305657
$vals[12,10] = 0.1093737608697887
$vals[12,11] = 0.1093737608697887
$vals[12,12] = 1.068379161376667
$vals[12,13] = 9.61541245239
$vals[12,14] = 0.01061365427691832
$vals[12,15] = 0.01061365427691832
$vals[13,0] = 1
$vals[13,1] = 0.3333333333333333
$vals[13,2] = 0.17509
$vals[13,3] = 0.52527
$vals[13,4] = 0.09704022420472538
$vals[13,5] = 0.09704022420472537
$vals[13,6] = 3
$vals[13,7] = 1
$vals[13,8] = 29.178763
$vals[13,9] = 87.53628900000001
$vals[13,10] = 0.5230171820937495
$vals[13,11] = 0.5230171820937495
$vals[13,12] = 5.108909613670001
$vals[13,13] = 45.98018652303001
$vals[13,14] = 0.05075370461330113
$vals[13,15] = 0.05075370461330112
$vals[14,0] = 1
$vals[14,1] = 0.3333333333333333
$vals[14,2] = 0.17509
$vals[14,3] = 0.52527
$vals[14,4] = 0.09704022420472538
$vals[14,5] = 0.09704022420472537
$vals[14,6] = 2
$vals[14,7] = 0.6666666666666666
$vals[14,8] = 0.146644
$vals[14,9] = 0.439932
$vals[14,10] = 0.002628532664354407
$vals[14,11] = 0.002628532664354407
$vals[14,12] = 0.02567589796
$vals[14,13] = 0.23108308164
$vals[14,14] = 0.0002550733990783958
$vals[14,15] = 0.0002550733990783958
$vals[15,0] = 1
$vals[15,1] = 0.3333333333333333
$vals[15,2] = 0.17509
$vals[15,3] = 0.52527
$vals[15,4] = 0.09704022420472538
$vals[15,5] = 0.09704022420472537
$vals[15,6] = 3
$vals[15,7] = 1
$vals[15,8] = 15.02284966666667
$vals[15,9] = 45.068549
$vals[15,10] = 0.2692783275177917
$vals[15,11] = 0.2692783275177917
$vals[15,12] = 2.630350748136666
$vals[15,13] = 23.67315673323
$vals[15,14] = 0.02613082927579998
$vals[15,15] = 0.02613082927579998
$vals[16,0] = 1
$vals[16,1] = 0.3333333333333333
$vals[16,2] = 0.17509
$vals[16,3] = 0.52527
$vals[16,4] = 0.09704022420472538
$vals[16,5] = 0.09704022420472537
$vals[16,6] = 3
$vals[16,7] = 1
$vals[16,8] = 5.288900666666667
$vals[16,9] = 15.866702
$vals[16,10] = 0.09480134312252211
$vals[16,11] = 0.09480134312252211
$vals[16,12] = 0.9260336177266667
$vals[16,13] = 8.334302559540001
$vals[16,14] = 0.009199543591518648
$vals[16,15] = 0.009199543591518646
$vals[17,0] = 1
$vals[17,1] = 0.3333333333333333
$vals[17,2] = 0.17509
$vals[17,3] = 0.52527
$vals[17,4] = 0.09704022420472538
$vals[17,5] = 0.09704022420472537
$vals[17,6] = 1
$vals[17,7] = 0.3333333333333333
$vals[17,8] = 0.050258
$vals[17,9] = 0.150774
$vals[17,10] = 0.0009008537317934847
$vals[17,11] = 0.0009008537317934848
$vals[17,12] = 0.00879967322
$vals[17,13] = 0.07919705898
$vals[17,14] = 0.0000874190481089033
$vals[17,15] = 0.0000874190481089033

$ws.Range("E2:T19").Value = $vals
Write-Output "applied"
